# Improve test file for detecting TimeSpan type cell.
# Add predefined (built-in) time formats 18-20 and 45-47 to the test sheet,
# in addition to the already-present format 21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rawFormat = $ws.Range("C3").NumberFormat   # "0.000000" -- the "Raw Value" column format

# --- Step 1: relocate the existing examples (originally rows 3-6) down to
# their new final location (rows 9-12), preserving their exact text/value/
# format so they keep the same shared-string / style identity.
$oldRows = 3,4,5,6
$newRows = 9,10,11,12

for ($i = 0; $i -lt $oldRows.Length; $i++) {
    $src = $oldRows[$i]
    $dst = $newRows[$i]

    $ws.Range("A$dst").Value = $ws.Range("A$src").Value()

    $ws.Range("B$dst").Value = $ws.Range("B$src").Value()
    $ws.Range("B$dst").NumberFormat = $ws.Range("B$src").NumberFormat
}

# Drop the old shared formula group living in C3:C6 (it will be overwritten
# by new content below) and rebuild it at its new location C9:C12 as a
# single shared formula, same as the original.
$ws.Range("C3:C6").ClearContents()
$ws.Range("C9:C12").Formula = "=B9"
$ws.Range("C9:C12").NumberFormat = $rawFormat

# --- Step 2: fill rows 2-8 with the new predefined-format examples.
# Row/format assignment (ascending row order controls the creation order of
# the new cell styles 18,19,20,45,46,47; format 21 reuses the existing style):
$byRow = @(
    @{ Row = 2; Format = "h:mm AM/PM" },
    @{ Row = 3; Format = "h:mm:ss AM/PM" },
    @{ Row = 4; Format = "h:mm" },
    @{ Row = 5; Format = "h:mm:ss" },
    @{ Row = 6; Format = "mm:ss" },
    @{ Row = 7; Format = "[h]:mm:ss" },
    @{ Row = 8; Format = "mm:ss.0" }
)
foreach ($item in $byRow) {
    $r = $item.Row
    $ws.Range("B$r").Value = 0.56313888888888886
    $ws.Range("B$r").NumberFormat = $item.Format
    $ws.Range("C$r").Formula = "=B$r"
    $ws.Range("C$r").NumberFormat = $rawFormat
}

# The label strings (column A) are introduced in this exact order so the
# shared-string table appends the new unique strings in this sequence.
$ws.Range("A5").Value = "Time fmt 21: h:mm:ss"
$ws.Range("A4").Value = "Time fmt 20: H:mm"
$ws.Range("A3").Value = "Time fmt 19: h:mm:ss tt"
$ws.Range("A2").Value = "Time fmt 18: h:mm tt"
$ws.Range("A6").Value = "Time fmt 45: mm:ss"
$ws.Range("A7").Value = "Time fmt 46: [h]:mm:ss"
$ws.Range("A8").Value = "Time fmt 47: mm:ss.0"

$ws.Range("C11").Select()

Write-Output "done"
